$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @("espino", "paolo", 502179, "espip001", "espinpa01", 8246, 2017, 2021, 29989),
    @("stroman", "marcus", 573186, "strom001", "stromma01", 13431, 2014, 2021, 32815),
    @("pineda", "michael", 501381, "pinem001", "pinedmi01", 5372, 2011, 2021, 30937),
    @("foltynewicz", "mike", 592314, "foltm001", "foltymi01", 10811, 2014, 2021, 31819),
    @("garcia", "luis", 671277, "garcl006", "garcilu04", -1, 2020, 2021, 4684365),
    @("ohtani", "shohei", 660271, "ohtas001", "ohtansh01", 19755, 2018, 2021, 39832),
    @("smith", "caleb", 592761, "smitc006", "smithca03", 14875, 2017, 2021, 36081),
    @("pivetta", "nick", 601713, "piven001", "pivetni01", 15454, 2017, 2021, 36071)
)

$startRow = 137
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
}
